$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -6.430699999999997
$ws.Range("E4").Value = 13.5295

$ws.Range("E5").Value = 13.553

$ws.Range("D6").Value = -7.902

$ws.Range("D7").Value = -7.672299999999995

$ws.Range("E8").Value = 14.2733

$ws.Range("D16").Value = -7.766899999999995
$ws.Range("E16").Value = 14.5712

$ws.Range("D20").Value = -8.025699999999995

$ws.Range("E22").Value = 13.46339999999999
